# Migration to Automation-Org/TestCases-maintenance/WIP-RMA TestCases
#
# The "Routing Master" sheet's sample row (row 2) references a specific
# Engineering Item record (its display name in column B, and its Salesforce
# record Id in column D). Re-point it at the newly migrated record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Routing Master")

$ws.Range("B2").Value = "Pro-PEItem-QQSVF"
$ws.Range("D2").Value = "a2S1K000002TIroUAG"

# Columns B and D are best-fit to their contents; refresh the fit now that
# the displayed text has changed.
$ws.Columns("B:B").AutoFit()
$ws.Columns("D:D").AutoFit()
